# Swap the two worksheets' identities:
#   - the sheet currently named "hotel_info" becomes "review_info"
#     (keeping only its header row, 25 columns: STR + the review_info fields)
#   - the sheet currently named "review_info" becomes "hotel_info"
#     (keeping its header row + single data row, now with a new "State"
#     column inserted right after "Hotel_Name")

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item(1)   # physically sheet1.xml, currently "hotel_info"
$wsReview = $wb.Worksheets.Item(2)   # physically sheet2.xml, currently "review_info"

# --- rename, via a temporary name so the two don't collide ---
$wsHotel.Name  = "__tmp_swap__"
$wsReview.Name = "hotel_info"
$wsHotel.Name  = "review_info"

# $wsHotel  is now named "review_info"
# $wsReview is now named "hotel_info"

# --- rebuild the (now) "review_info" sheet: header row only ---
$wsHotel.Cells.ClearContents()

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $wsHotel.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- rebuild the (now) "hotel_info" sheet: header row + data row, with new State column ---
$wsReview.Cells.ClearContents()

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $wsReview.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$wsReview.Cells.Item(2, 1).Value = 1672
$wsReview.Cells.Item(2, 2).Value = "Sheraton Hotel New Orleans"
$wsReview.Cells.Item(2, 3).Value = "Louisiana"
$wsReview.Cells.Item(2, 4).Value = "New Orleans"
$wsReview.Cells.Item(2, 5).Value = 70130
$wsReview.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d89130-Reviews-Sheraton_New_Orleans_Hotel-New_Orleans_Louisiana.html"
$wsReview.Cells.Item(2, 7).Value = "Sheraton New Orleans Hotel"

# English_Reviews_num / Local_Rank / Total_Reviews_num are stored as text
# (numeric-looking strings) in the source workbook, not as numbers - force
# text storage, then drop back to the Normal style so no stray number
# format sticks around on the cell.
$wsReview.Range("H2").NumberFormat = "@"
$wsReview.Cells.Item(2, 8).Value = "2384"
$wsReview.Range("H2").Style = "Normal"

$wsReview.Range("I2").NumberFormat = "@"
$wsReview.Cells.Item(2, 9).Value = "99"
$wsReview.Range("I2").Style = "Normal"

$wsReview.Range("J2").NumberFormat = "@"
$wsReview.Cells.Item(2, 10).Value = "2508"
$wsReview.Range("J2").Style = "Normal"
